# Updated symbol list on Sun Jan 15 09:44:59 UTC 2023 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) figures on the
# "Coin" table for the rows whose quotes moved since the previous snapshot.
# Values are written as literal text (matching the sheet's existing
# inline-string cell contents, e.g. "294.69" / "-4.11%") rather than as
# numbers/percentages, so we briefly force a Text number format while
# assigning the value and then restore the cell to the workbook's default
# "Normal" style (this sheet does not use custom number formats on these
# cells) to avoid leaving any formatting residue behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue([string]$ref, [string]$newValue) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $newValue
    $cell.Style = "Normal"
}


Set-TextValue "D2" "294.69"
Set-TextValue "E2" "-4.11%"
Set-TextValue "D3" "31.25"
Set-TextValue "E3" "-1.83%"
Set-TextValue "D4" "5.107"
Set-TextValue "E4" "-3.14%"
Set-TextValue "D5" "0.07357"
Set-TextValue "E5" "0.52%"
Set-TextValue "D6" "7.672"
Set-TextValue "E6" "-2.28%"
Set-TextValue "D7" "3.756"
Set-TextValue "E7" "0.05%"
Set-TextValue "D8" "1.642"
Set-TextValue "E8" "10.45%"
Set-TextValue "D9" "0.9203"
Set-TextValue "E9" "1.58%"
Set-TextValue "D10" "0.1668"
Set-TextValue "E10" "-1.08%"
Set-TextValue "D11" "0.07042"
Set-TextValue "E11" "-6.19%"
Set-TextValue "D12" "0.07972"
Set-TextValue "E12" "-0.95%"
Set-TextValue "D13" "0.02999"
Set-TextValue "E13" "-0.15%"
Set-TextValue "D14" "0.09891"
Set-TextValue "E14" "-0.47%"
Set-TextValue "D15" "0.001489"
Set-TextValue "E15" "-1.40%"
Set-TextValue "D16" "0.006160"
Set-TextValue "E16" "-2.27%"
Set-TextValue "D17" "3.452"
Set-TextValue "E17" "-0.84%"
Set-TextValue "E18" "-0.14%"
Set-TextValue "D19" "0.3278"
Set-TextValue "E19" "-1.45%"
Set-TextValue "D20" "0.1334"
Set-TextValue "E20" "-1.07%"
Set-TextValue "D21" "4.561"
Set-TextValue "E21" "5.44%"
Set-TextValue "D22" "0.04626"
Set-TextValue "E22" "1.05%"
Set-TextValue "E23" "-5.23%"
Set-TextValue "D24" "0.001216"
Set-TextValue "E24" "-0.99%"
Set-TextValue "D25" "0.004415"
Set-TextValue "E25" "-0.50%"
Set-TextValue "D26" "0.0001301"
Set-TextValue "E26" "-0.51%"
Set-TextValue "E27" "7.14%"
Set-TextValue "E39" "1.79%"
Set-TextValue "D40" "0.04407"
Set-TextValue "E40" "-2.08%"
Set-TextValue "D41" "0.007114"
Set-TextValue "E41" "-1.78%"
Set-TextValue "D42" "0.1326"
Set-TextValue "E42" "-1.58%"
Set-TextValue "D43" "0.002112"
Set-TextValue "E43" "-7.10%"
Set-TextValue "D44" "0.01100"
Set-TextValue "E44" "-22.05%"
Set-TextValue "D45" "0.00005991"
Set-TextValue "E45" "-1.27%"
Set-TextValue "D46" "0.7116"
Set-TextValue "E46" "-62.40%"
Set-TextValue "E47" "-15.85%"
